# Applies the "Subo la base corregida en excel" correction:
#  - Fixes three MySQL-for-Excel Data Model defined names (dropping a stray
#    trailing "1" so they match their connections' range names).
#  - Resizes the affected ListObjects (Tabla6, Tabla13, Tabla16, Tabla18) so
#    their ranges track the row that's about to become their new header.
#  - Removes the stray "id_cordenadas" label that had been left in D7 and
#    shifts the rest of column D (the several small lookup tables below it)
#    up by one row to close the gap.
#  - Restores the last saved selection to E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Defined names: ...articulos1 -> ...articulos, ...tipocontenido1 ->
#    ...tipocontenido, ...usuarios1 -> ...usuarios (the "...articulos11"
#    name is left untouched).
# ---------------------------------------------------------------------
$wb.Names.Item("_xlcn.ModelConnection_For_Libro1zigma.articulos1").Name = "_xlcn.ModelConnection_For_Libro1zigma.articulos"
$wb.Names.Item("_xlcn.ModelConnection_For_Libro1zigma.tipocontenido1").Name = "_xlcn.ModelConnection_For_Libro1zigma.tipocontenido"
$wb.Names.Item("_xlcn.ModelConnection_For_Libro1zigma.usuarios1").Name = "_xlcn.ModelConnection_For_Libro1zigma.usuarios"

# ---------------------------------------------------------------------
# 2. Resize the tables that live in column D to their new ranges BEFORE
#    touching any cell values, so each table picks up its (about to be
#    relocated) header text from the still-correct, pre-shift cells.
# ---------------------------------------------------------------------
$ws.ListObjects("Tabla6").Resize($ws.Range("D2:D6"))
$ws.ListObjects("Tabla13").Resize($ws.Range("D9:D13"))
$ws.ListObjects("Tabla16").Resize($ws.Range("D15:D19"))
$ws.ListObjects("Tabla18").Resize($ws.Range("D22:D26"))

# ---------------------------------------------------------------------
# 3. Column D, rows 7-27: delete D7 ("id_cordenadas") and shift the cells
#    below it up by one row. Read all the old values first (top to
#    bottom) so the write-back pass below never clobbers data it still
#    needs to read.
# ---------------------------------------------------------------------
$colD = 4
$firstRow = 7
$lastRow = 27

$shiftedValues = @()
for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $shiftedValues += ,$ws.Cells.Item($r, $colD).Text
}

for ($r = $firstRow; $r -le ($lastRow - 1); $r++) {
    $v = $shiftedValues[$r - $firstRow]
    if ($v -eq "") {
        $ws.Cells.Item($r, $colD).ClearContents() | Out-Null
    } else {
        $ws.Cells.Item($r, $colD).Value = $v
    }
}
$ws.Cells.Item($lastRow, $colD).ClearContents() | Out-Null

# ---------------------------------------------------------------------
# 4. Restore the saved selection (E6).
# ---------------------------------------------------------------------
$ws.Range("E6").Select() | Out-Null
